$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 78, shifting all
# existing data (rows 78-182) down to rows 82-186.
$ws.Range("A78:A81").EntireRow.Insert()

# Populate the 4 newly inserted rows with the new weekly price records.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg,
# Q Kg o Unidades, R Clasificación

$newRows = @(
    @{ Row=78; Calidad="Cuarta";  Volumen=2000; PMin=1500; PMax=1500; PProm=1500; PKg=1500 },
    @{ Row=79; Calidad="Primera"; Volumen=2000; PMin=3000; PMax=3000; PProm=3000; PKg=3000 },
    @{ Row=80; Calidad="Segunda"; Volumen=2000; PMin=2500; PMax=2500; PProm=2500; PKg=2500 },
    @{ Row=81; Calidad="Tercera"; Volumen=2000; PMin=2000; PMax=2000; PProm=2000; PKg=2000 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value2 = 4
    $ws.Cells.Item($r, 2).Value2 = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($r, 3).Value2 = "Los Lagos"
    $ws.Cells.Item($r, 4).Value2 = 44571
    $ws.Cells.Item($r, 5).Value2 = 10
    $ws.Cells.Item($r, 6).Value2 = 100112028
    $ws.Cells.Item($r, 7).Value2 = "Sandia"
    $ws.Cells.Item($r, 8).Value2 = "Sin especificar"
    $ws.Cells.Item($r, 9).Value2 = $rec.Calidad
    $ws.Cells.Item($r, 10).Value2 = $rec.Volumen
    $ws.Cells.Item($r, 11).Value2 = $rec.PMin
    $ws.Cells.Item($r, 12).Value2 = $rec.PMax
    $ws.Cells.Item($r, 13).Value2 = $rec.PProm
    $ws.Cells.Item($r, 14).Value2 = "`$/unidad"
    $ws.Cells.Item($r, 15).Value2 = "Región del Maule"
    $ws.Cells.Item($r, 16).Value2 = $rec.PKg
    $ws.Cells.Item($r, 17).Value2 = 1
    $ws.Cells.Item($r, 18).Value2 = "Hortaliza"
}
